$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'290.54"
$ws.Range("E2").Value = "'1.64%"
$ws.Range("D3").Value = "'30.77"
$ws.Range("E3").Value = "'-2.48%"
$ws.Range("D4").Value = "'4.869"
$ws.Range("E4").Value = "'-2.18%"
$ws.Range("D5").Value = "'0.07248"
$ws.Range("E5").Value = "'-0.37%"
$ws.Range("D6").Value = "'2.433"
$ws.Range("E6").Value = "'36.69%"
$ws.Range("D7").Value = "'7.644"
$ws.Range("E7").Value = "'0.22%"
$ws.Range("D8").Value = "'3.705"
$ws.Range("E8").Value = "'-0.32%"
$ws.Range("D9").Value = "'0.8963"
$ws.Range("E9").Value = "'-1.58%"
$ws.Range("D10").Value = "'0.1669"
$ws.Range("E10").Value = "'2.15%"
$ws.Range("D11").Value = "'0.08053"
$ws.Range("E11").Value = "'7.58%"
$ws.Range("D12").Value = "'0.08155"
$ws.Range("E12").Value = "'-0.23%"
$ws.Range("D13").Value = "'0.03072"
$ws.Range("E13").Value = "'2.21%"
$ws.Range("D14").Value = "'0.1003"
$ws.Range("E14").Value = "'0.31%"
$ws.Range("D15").Value = "'0.001497"
$ws.Range("E15").Value = "'-0.29%"
$ws.Range("D16").Value = "'0.005739"
$ws.Range("E16").Value = "'-0.45%"
$ws.Range("D17").Value = "'3.485"
$ws.Range("E17").Value = "'0.51%"
$ws.Range("D18").Value = "'2.076"
$ws.Range("E18").Value = "'-2.13%"
$ws.Range("D19").Value = "'0.3314"
$ws.Range("E19").Value = "'1.67%"
$ws.Range("D20").Value = "'0.1287"
$ws.Range("E20").Value = "'-0.56%"
$ws.Range("D21").Value = "'3.969"
$ws.Range("E21").Value = "'-9.48%"
$ws.Range("E23").Value = "'0.13%"
$ws.Range("D24").Value = "'0.001212"
$ws.Range("E24").Value = "'-2.07%"
$ws.Range("D25").Value = "'0.004406"
$ws.Range("E25").Value = "'10.11%"
$ws.Range("D26").Value = "'0.0001301"
$ws.Range("E26").Value = "'2.92%"
$ws.Range("D39").Value = "'0.01588"
$ws.Range("E39").Value = "'-1.12%"
$ws.Range("D40").Value = "'0.04376"
$ws.Range("E40").Value = "'1.06%"
$ws.Range("D41").Value = "'0.007251"
$ws.Range("E41").Value = "'-2.65%"
$ws.Range("D42").Value = "'0.01001"
$ws.Range("D43").Value = "'0.1313"
$ws.Range("E43").Value = "'-0.15%"
$ws.Range("D44").Value = "'0.002025"
$ws.Range("E44").Value = "'-6.51%"
$ws.Range("D45").Value = "'0.009170"
$ws.Range("E45").Value = "'-14.58%"
$ws.Range("D46").Value = "'0.00005711"
$ws.Range("E46").Value = "'-6.40%"
$ws.Range("E47").Value = "'-1.04%"
$ws.Range("D48").Value = "'2.241"
$ws.Range("E48").Value = "'18.18%"
$ws.Range("E49").Value = "'-4.35%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-1.04%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-1.04%"
